# edit.ps1
# Applies the Nalco aluminium price update:
#   - Inserts a new "top" record (date 23-12-2025) that duplicates the
#     previous top record's Description/Product Code/Basic Price/Circular
#     Date/Circular Link (only the Date column differs), pushing every
#     existing data row down by one row.
#   - The very last existing data row therefore ends up duplicated at the
#     new bottom of the table (row 140), matching the source diff.
#   - Dimension grows from A1:F139 to A1:F140 and a new hyperlink is added
#     for the new bottom row's Circular Link cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- discover current extent of the table -----------------------------
$usedRange  = $ws.UsedRange
$headerRow  = $usedRange.Row()                      # normally 1
$lastRow    = $headerRow + $usedRange.Rows.Count() - 1   # last data row (139)
$firstData  = $headerRow + 1                         # first data row (2)
$lastCol    = 6                                      # columns A..F

# --- stash the "normal data row" cell style somewhere far away so we can
#     restore it later on cells whose formatting gets clobbered by
#     Hyperlinks.Add (it forces the built-in "Hyperlink" style) ----------
$stashCell = $ws.Cells.Item($lastRow + 500, $lastCol + 10)
$ws.Cells.Item($firstData, $lastCol).Copy($stashCell)

# --- shift every existing data row down by one row, bottom-up so we never
#     overwrite a row before it has been copied -------------------------
for ($r = $lastRow; $r -ge $firstData; $r--) {
    $destRow    = $r + 1
    $srcRange   = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))
    $destRange  = $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, $lastCol))
    $srcRange.Copy($destRange)
}

$newLastRow = $lastRow + 1

# --- write the new top row's date (rest of the row already matches the
#     old top row, which Copy() just duplicated one row down) -----------
$ws.Cells.Item($firstData, 1).Value2 = "23-12-2025"

# --- rebuild every hyperlink in the Circular Link column (F) from the
#     text now sitting in each cell; the old Hyperlinks collection does
#     not auto-shift when rows are copied, so start clean ---------------
$ws.Hyperlinks.Delete()
for ($r = $firstData; $r -le $newLastRow; $r++) {
    $linkCell = $ws.Cells.Item($r, $lastCol)
    $url = $linkCell.Value2()
    $ws.Hyperlinks.Add($linkCell, $url) | Out-Null
}

# --- restore the original (non-hyperlink) cell style on the whole Link
#     column, since Hyperlinks.Add() force-applied the "Hyperlink" style
$stashCell.Copy()
$linkColRange = $ws.Range($ws.Cells.Item($firstData, $lastCol), $ws.Cells.Item($newLastRow, $lastCol))
$linkColRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- cleanup ------------------------------------------------------------
$stashCell.Clear() | Out-Null
$excel.CutCopyMode = 0

"Updated sheet: rows $firstData..$newLastRow (was $firstData..$lastRow); dimension A1:F$newLastRow"
